# Move the "to" address (column C) for every shipment row onto a single
# consolidated postcode ("M16 0RA"), collapsing the two near-duplicate
# "M41 4LR" / "M414LR" destination strings that used to live side by side
# in the shared-string table.
#
# Rows 12, 31, 127 and 166 are intentional blank "gap" rows (no B/C data);
# their column C cell used a stray numeric style left over from the old
# data - bring it into line with the rest of column C (text style) now
# that the column is uniformly textual.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 272
$gapRows = @(12, 31, 127, 166)
$destination = "M16 0RA"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    if ($gapRows -contains $r) { continue }
    $ws.Cells.Item($r, 3).Value = $destination
}

foreach ($r in $gapRows) {
    $ws.Cells.Item($r, 3).NumberFormat = "@"
}
